$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to remain text so values such as
# "1.00" or "68.455.07" are not re-interpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "68.455.07"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "3.855.28"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "521.23"
$ws.Range("E5").Value = "  +6.81%  "
$ws.Range("D6").Value = "140.57"
$ws.Range("E6").Value = "  -3.94%  "
$ws.Range("E7").Value = "  -2.47%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("E10").Value = "  -5.93%  "
$ws.Range("D11").Value = "0.0000318"
$ws.Range("E11").Value = "  -7.85%  "
$ws.Range("D12").Value = "41.43"
$ws.Range("E12").Value = "  -3.51%  "
$ws.Range("D13").Value = "10.32"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "4.466.66"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "21.24"
$ws.Range("E15").Value = "  +6.59%  "
$ws.Range("D16").Value = "3.849.98"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").Value = "14.12"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "68.484.77"
$ws.Range("D21").Value = "415.56"
$ws.Range("E21").Value = "  -4.94%  "
$ws.Range("D22").Value = "3.48"
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "13.95"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").Value = "86.61"
$ws.Range("E24").Value = "  -3.02%  "
$ws.Range("E25").Value = "  +6.20%  "
$ws.Range("D26").Value = "11.50"
$ws.Range("E26").Value = "  -6.74%  "
$ws.Range("D27").Value = "10.53"
$ws.Range("E27").Value = "  -5.26%  "
$ws.Range("D28").Value = "35.43"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "13.18"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").Value = "677.50"
$ws.Range("E31").Value = "  -6.34%  "
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("D33").Value = "6.61"
$ws.Range("E33").Value = "  +9.07%  "
$ws.Range("D34").Value = "65.89"
$ws.Range("E34").Value = "  +6.79%  "
$ws.Range("D35").Value = "0.450"
$ws.Range("E35").Value = "  -6.08%  "
$ws.Range("D36").Value = "39.55"
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("D37").Value = "3.49"
$ws.Range("E37").Value = "  +13.96%  "
$ws.Range("D38").Value = "0.0₃0827"
$ws.Range("E38").Value = "  -7.95%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E43").Value = "  +2.66%  "
$ws.Range("E44").Value = "  -5.38%  "
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "0.139"
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "143.29"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D49").Value = "3.25"
$ws.Range("E49").Value = "  -3.93%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "0.000264"
$ws.Range("E50").Value = "  +11.11%  "
$ws.Range("D51").Value = "0.0₆0336"
$ws.Range("E51").Value = "  -6.01%  "
